# Auto-generated edit script applying the Phantom_Profits diff
# Updates numeric cell values across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8081.5
$ws.Range("I32").Value = 8387
$ws.Range("K32").Value = 8387
$ws.Range("M32").Value = -8061
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H86").Value = 6125.75
$ws.Range("I86").Value = 6000
$ws.Range("J86").Value = 6251.5
$ws.Range("K86").Value = 6000
$ws.Range("L86").Value = 6251.5
$ws.Range("M86").Value = -4877
$ws.Range("N86").Value = -8497.5
$ws.Range("H89").Value = 6125.75
$ws.Range("I89").Value = 6000
$ws.Range("J89").Value = 6251.5
$ws.Range("K89").Value = 30000
$ws.Range("L89").Value = 31257.5
$ws.Range("M89").Value = -24384
$ws.Range("N89").Value = -42489.5
$ws.Range("H100").Value = 2206.7058
$ws.Range("I100").Value = 2065.7144
$ws.Range("K100").Value = 2065.7144
$ws.Range("M100").Value = -1524.7144
$ws.Range("H127").Value = 1071.2
$ws.Range("I127").Value = 964.25
$ws.Range("J127").Value = 1499
$ws.Range("K127").Value = 2892.75
$ws.Range("L127").Value = 4497
$ws.Range("M127").Value = 2067.25
$ws.Range("N127").Value = -14417
$ws.Range("H132").Value = 6970
$ws.Range("I132").Value = 7180.875
$ws.Range("J132").Value = 3596
$ws.Range("K132").Value = 21542.625
$ws.Range("L132").Value = 10788
$ws.Range("M132").Value = -19012.625
$ws.Range("N132").Value = -15848

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3205.6538
$ws.Range("I32").Value = 2933.92
$ws.Range("K32").Value = 2933.92
$ws.Range("M32").Value = -2646.92
$ws.Range("H61").Value = 2559.8
$ws.Range("I61").Value = 2559.8
$ws.Range("K61").Value = 2559.8
$ws.Range("M61").Value = -2347.8
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H136").Value = 2559.8
$ws.Range("I136").Value = 2559.8
$ws.Range("K136").Value = 7679.400000000001
$ws.Range("M136").Value = -5129.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3172.2727
$ws.Range("I20").Value = 2089.1428
$ws.Range("J20").Value = 5067.75
$ws.Range("K20").Value = 2089.1428
$ws.Range("L20").Value = 5067.75
$ws.Range("M20").Value = -1842.1428
$ws.Range("N20").Value = -5561.75
$ws.Range("H86").Value = 3453.8
$ws.Range("I86").Value = 3453.8
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3453.8
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3453.8
$ws.Range("I89").Value = 3453.8
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17269
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H105").Value = 3141
$ws.Range("I105").Value = 2405.8262
$ws.Range("K105").Value = 2405.8262
$ws.Range("M105").Value = -658.8262
$ws.Range("H107").Value = 834
$ws.Range("I107").Value = 622.0909
$ws.Range("K107").Value = 622.0909
$ws.Range("M107").Value = 1297.9091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 100000
$ws.Range("J9").Value = 100000
$ws.Range("L9").Value = 100000
$ws.Range("N9").Value = -100336
$ws.Range("H35").Value = 880
$ws.Range("I35").Value = 880
$ws.Range("K35").Value = 880
$ws.Range("M35").Value = -586
$ws.Range("H74").Value = 18000
$ws.Range("J74").Value = 18000
$ws.Range("L74").Value = 18000
$ws.Range("N74").Value = -19748
$ws.Range("H77").Value = 18000
$ws.Range("J77").Value = 18000
$ws.Range("L77").Value = 54000
$ws.Range("N77").Value = -62736
$ws.Range("H88").Value = 9920
$ws.Range("J88").Value = 9920
$ws.Range("L88").Value = 9920
$ws.Range("N88").Value = -10732
$ws.Range("H91").Value = 9920
$ws.Range("J91").Value = 9920
$ws.Range("L91").Value = 9920
$ws.Range("N91").Value = -12728
$ws.Range("H94").Value = 1312.5
$ws.Range("J94").Value = 700
$ws.Range("L94").Value = 700
$ws.Range("N94").Value = -1602
$ws.Range("H99").Value = 2998.9443
$ws.Range("I99").Value = 2891.6667
$ws.Range("J99").Value = 3213.5
$ws.Range("K99").Value = 2891.6667
$ws.Range("L99").Value = 3213.5
$ws.Range("M99").Value = -1393.6667
$ws.Range("N99").Value = -6209.5
$ws.Range("H126").Value = 2998.9443
$ws.Range("I126").Value = 2891.6667
$ws.Range("J126").Value = 3213.5
$ws.Range("K126").Value = 8675.000100000001
$ws.Range("L126").Value = 9640.5
$ws.Range("M126").Value = -6205.000100000001
$ws.Range("N126").Value = -14580.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H4").Value = 1683887.5
$ws.Range("I4").Value = 18233.883
$ws.Range("K4").Value = 54701.649
$ws.Range("M4").Value = -54589.649
$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 60000
$ws.Range("N88").Value = -60856
$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 60000
$ws.Range("N91").Value = -62964
$ws.Range("H132").Value = 1181.6923
$ws.Range("I132").Value = 1295.8889
$ws.Range("J132").Value = 924.75
$ws.Range("K132").Value = 11663.0001
$ws.Range("L132").Value = 8322.75
$ws.Range("M132").Value = -9133.000099999999
$ws.Range("N132").Value = -13382.75
$ws.Range("H133").Value = 17535.75
$ws.Range("I133").Value = 9809.666999999999
$ws.Range("K133").Value = 29429.001
$ws.Range("M133").Value = -24369.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 999
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H16").Value = 999
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 29900
$ws.Range("J46").Value = 29900
$ws.Range("L46").Value = 29900
$ws.Range("N46").Value = -30212
$ws.Range("H70").Value = 6937
$ws.Range("I70").Value = 6947.75
$ws.Range("K70").Value = 6947.75
$ws.Range("M70").Value = -6677.75
$ws.Range("H73").Value = 6937
$ws.Range("I73").Value = 6947.75
$ws.Range("K73").Value = 6947.75
$ws.Range("M73").Value = -6011.75
$ws.Range("H128").Value = 47296.6
$ws.Range("I128").Value = 40500
$ws.Range("J128").Value = 48995.75
$ws.Range("K128").Value = 40500
$ws.Range("L128").Value = 48995.75
$ws.Range("M128").Value = -35520
$ws.Range("N128").Value = -58955.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 546.5
$ws.Range("I9").Value = 546.5
$ws.Range("K9").Value = 546.5
$ws.Range("M9").Value = -322.5
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9774
$ws.Range("H30").Value = 1108.25
$ws.Range("I30").Value = 1372
$ws.Range("J30").Value = 317
$ws.Range("K30").Value = 1372
$ws.Range("L30").Value = 317
$ws.Range("M30").Value = -1264
$ws.Range("N30").Value = -533

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 66228
$ws.Range("I45").Value = 39997.5
$ws.Range("K45").Value = 39997.5
$ws.Range("M45").Value = -39506.5
